# Sync attendance_reports, modules_schedules, and assets from main repo - 2026-01-05 16:21:07
#
# The "Recorded By" column (G) on the "Session Analysis Results" sheet lists the
# users who recorded a session, separated by ", ". This edit reverses the order
# of that comma-separated list for each row - except for rows whose value is the
# literal "admin@admin.com, System" (left untouched) and rows that only contain a
# single name/email (nothing to reorder).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -isnot [string]) { continue }
    if ($val -eq "admin@admin.com, System") { continue }
    if ($val.IndexOf(",") -lt 0) { continue }

    $parts = $val.Split(",")
    for ($i = 0; $i -lt $parts.Length; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }
    $n = $parts.Length
    $reversed = $parts[($n - 1)..0]
    $newVal = [string]::Join(", ", $reversed)

    if ($newVal -ne $val) {
        $cell.Value = $newVal
    }
}
